# Auto-generated: apply Leviathan_Profits market-data refresh
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 5000.3335
$ws.Range("I62").Value = 4469.091
$ws.Range("K62").Value = 4469.091
$ws.Range("M62").Value = -3845.091
$ws.Range("H65").Value = 5000.3335
$ws.Range("I65").Value = 4469.091
$ws.Range("K65").Value = 22345.455
$ws.Range("M65").Value = -19225.455
$ws.Range("H98").Value = 2453.25
$ws.Range("I98").Value = 1412.8572
$ws.Range("K98").Value = 1412.8572
$ws.Range("M98").Value = 85.14280000000008
$ws.Range("H122").Value = 2453.25
$ws.Range("I122").Value = 1412.8572
$ws.Range("K122").Value = 4238.571599999999
$ws.Range("M122").Value = -1788.571599999999
$ws.Range("H133").Value = 96000
$ws.Range("J133").Value = 96000
$ws.Range("L133").Value = 96000
$ws.Range("N133").Value = -106120
$ws.Range("H138").Value = 2291.5757
$ws.Range("I138").Value = 2557.7273
$ws.Range("J138").Value = 2158.5
$ws.Range("K138").Value = 7673.1819
$ws.Range("L138").Value = 6475.5
$ws.Range("M138").Value = -2533.1819
$ws.Range("N138").Value = -16755.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3081.2307
$ws.Range("I32").Value = 3113.5134
$ws.Range("K32").Value = 3113.5134
$ws.Range("M32").Value = -2826.5134
$ws.Range("H45").Value = 5457.086
$ws.Range("I45").Value = 7424.737
$ws.Range("J45").Value = 3120.5
$ws.Range("K45").Value = 7424.737
$ws.Range("L45").Value = 3120.5
$ws.Range("M45").Value = -7047.737
$ws.Range("N45").Value = -3874.5
$ws.Range("H61").Value = 3339.3333
$ws.Range("I61").Value = 3339.3333
$ws.Range("K61").Value = 3339.3333
$ws.Range("M61").Value = -3127.3333
$ws.Range("H74").Value = 2906.9707
$ws.Range("I74").Value = 2697.2964
$ws.Range("K74").Value = 2697.2964
$ws.Range("M74").Value = -1823.2964
$ws.Range("H77").Value = 2906.9707
$ws.Range("I77").Value = 2697.2964
$ws.Range("K77").Value = 13486.482
$ws.Range("M77").Value = -9118.482
$ws.Range("H110").Value = 2715.4285
$ws.Range("I110").Value = 2715.4285
$ws.Range("K110").Value = 2715.4285
$ws.Range("M110").Value = -670.4285
$ws.Range("H122").Value = 2009.5385
$ws.Range("I122").Value = 2067.5454
$ws.Range("K122").Value = 6202.6362
$ws.Range("M122").Value = -3752.6362
$ws.Range("H136").Value = 3339.3333
$ws.Range("I136").Value = 3339.3333
$ws.Range("K136").Value = 10017.9999
$ws.Range("M136").Value = -7467.999899999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 4219.7
$ws.Range("J94").Value = 5199.6
$ws.Range("L94").Value = 5199.6
$ws.Range("N94").Value = -6101.6
$ws.Range("H134").Value = 1082.5
$ws.Range("I134").Value = 849.05
$ws.Range("J134").Value = 1860.6666
$ws.Range("K134").Value = 2547.15
$ws.Range("L134").Value = 5581.9998
$ws.Range("M134").Value = -12.14999999999964
$ws.Range("N134").Value = -10651.9998

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 3687.111
$ws.Range("I86").Value = 2990.8462
$ws.Range("K86").Value = 2990.8462
$ws.Range("M86").Value = -1867.8462
$ws.Range("H89").Value = 3687.111
$ws.Range("I89").Value = 2990.8462
$ws.Range("K89").Value = 14954.231
$ws.Range("M89").Value = -9338.231
$ws.Range("H132").Value = 4161.467
$ws.Range("I132").Value = 3878.6155
$ws.Range("K132").Value = 11635.8465
$ws.Range("M132").Value = -9105.8465
$ws.Range("H134").Value = 2546.1785
$ws.Range("I134").Value = 2307.5
$ws.Range("J134").Value = 3421.3333
$ws.Range("K134").Value = 6922.5
$ws.Range("L134").Value = 10263.9999
$ws.Range("M134").Value = -4387.5
$ws.Range("N134").Value = -15333.9999
$ws.Range("H141").Value = 439833
$ws.Range("J141").Value = 439833
$ws.Range("L141").Value = 439833
$ws.Range("N141").Value = -450193

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 1066.4445
$ws.Range("I8").Value = 1066.4445
$ws.Range("K8").Value = 3199.3335
$ws.Range("M8").Value = -3060.3335
$ws.Range("H34").Value = 359
$ws.Range("I34").Value = 152
$ws.Range("K34").Value = 456
$ws.Range("M34").Value = -372

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3066.0303
$ws.Range("I80").Value = 2719.7827
$ws.Range("J80").Value = 3862.4
$ws.Range("K80").Value = 2719.7827
$ws.Range("L80").Value = 3862.4
$ws.Range("M80").Value = -1721.7827
$ws.Range("N80").Value = -5858.4
$ws.Range("H83").Value = 3066.0303
$ws.Range("I83").Value = 2719.7827
$ws.Range("J83").Value = 3862.4
$ws.Range("K83").Value = 13598.9135
$ws.Range("L83").Value = 19312
$ws.Range("M83").Value = -8606.913500000001
$ws.Range("N83").Value = -29296
$ws.Range("H113").Value = 3131.4
$ws.Range("I113").Value = 3554.2222
$ws.Range("K113").Value = 3554.2222
$ws.Range("M113").Value = -1384.2222
$ws.Range("H126").Value = 4099.25
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 5100
$ws.Range("I68").Value = 5000
$ws.Range("J68").Value = 5200
$ws.Range("K68").Value = 5000
$ws.Range("L68").Value = 5200
$ws.Range("M68").Value = -4251
$ws.Range("N68").Value = -6698
$ws.Range("H71").Value = 5100
$ws.Range("I71").Value = 5000
$ws.Range("J71").Value = 5200
$ws.Range("K71").Value = 25000
$ws.Range("L71").Value = 26000
$ws.Range("M71").Value = -21256
$ws.Range("N71").Value = -33488
$ws.Range("H96").Value = 30000
$ws.Range("J96").Value = 30000
$ws.Range("L96").Value = 30000
$ws.Range("N96").Value = -35492
$ws.Range("H132").Value = 2673.3674
$ws.Range("I132").Value = 2197.9512
$ws.Range("J132").Value = 5109.875
$ws.Range("K132").Value = 6593.8536
$ws.Range("L132").Value = 15329.625
$ws.Range("M132").Value = -4063.8536
$ws.Range("N132").Value = -20389.625

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3024.5
$ws.Range("I81").Value = 958.1667
$ws.Range("K81").Value = 1916.3334
$ws.Range("M81").Value = -855.3334
$ws.Range("H84").Value = 3024.5
$ws.Range("I84").Value = 958.1667
$ws.Range("K84").Value = 9581.666999999999
$ws.Range("M84").Value = -4277.666999999999
$ws.Range("H107").Value = 13890229
$ws.Range("I107").Value = 1121.5834
$ws.Range("J107").Value = 41668444
$ws.Range("K107").Value = 3364.7502
$ws.Range("L107").Value = 125005332
$ws.Range("M107").Value = -1444.7502
$ws.Range("N107").Value = -125009172

Write-Host "Applied all profit/price updates to Leviathan_Profits workbook"
